# Add two new "long form" data tables to the workbook:
#   Median_HH_Income  -> race / med_hh_income
#   Poverty_HH        -> race / percentage
# and register them on the TOC sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Create the Median_HH_Income sheet (placed after "Education")
# ---------------------------------------------------------------
$education = $wb.Worksheets.Item("Education")
$medIncome = $wb.Worksheets.Add($null, $education)
$medIncome.Name = "Median_HH_Income"

$medIncome.Range("A1").Value = "race"
$medIncome.Range("B1").Value = "med_hh_income"
$medIncome.Range("A1:B1").Font.Bold = $true
$medIncome.Range("A1:B1").HorizontalAlignment = -4108

$medIncomeData = @(
    @("Asian/Pacific Islander", 93517.74),
    @("Black", 57918.45),
    @("Hispanic/Latino", 66943.16),
    @("Multiracial/Other", 86696.87),
    @("Native American", 65439.84),
    @("White", 96171.5),
    @("SCAG region", 79645.22)
)

$r = 2
foreach ($row in $medIncomeData) {
    $medIncome.Cells.Item($r, 1).Value = $row[0]
    $medIncome.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# ---------------------------------------------------------------
# 2. Create the Poverty_HH sheet (placed after "Median_HH_Income")
# ---------------------------------------------------------------
$poverty = $wb.Worksheets.Add($null, $medIncome)
$poverty.Name = "Poverty_HH"

$poverty.Range("A1").Value = "race"
$poverty.Range("B1").Value = "percentage"
$poverty.Range("A1:B1").Font.Bold = $true
$poverty.Range("A1:B1").HorizontalAlignment = -4108

$povertyData = @(
    @("Asian/Pacific Islander", 25.42),
    @("Black", 37.09),
    @("Hispanic/Latino", 38.39),
    @("Multiracial/Other", 24.82),
    @("Native American", 35.52),
    @("White", 19.95),
    @("SCAG", 28.89)
)

$r = 2
foreach ($row in $povertyData) {
    $poverty.Cells.Item($r, 1).Value = $row[0]
    $poverty.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# ---------------------------------------------------------------
# 3. Register both new sheets on the TOC sheet (sheet1)
# ---------------------------------------------------------------
$toc = $wb.Worksheets.Item("TOC")

$toc.Range("A22").Value = "Median_HH_Income"
$toc.Range("B22").Value = "Median Household Income by Race/ethincity and whole SCAG Region"

$toc.Range("A23").Value = "Poverty_HH"
$toc.Range("B23").Value = "Household Poverty (%) by Race/ethincity and whole SCAG Region"
